$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3303.4
$ws.Range("I69").Value = 3373.1667
$ws.Range("J69").Value = 3198.75
$ws.Range("K69").Value = 10119.5001
$ws.Range("L69").Value = 9596.25
$ws.Range("M69").Value = -9245.500100000001
$ws.Range("N69").Value = -11344.25
$ws.Range("H72").Value = 3303.4
$ws.Range("I72").Value = 3373.1667
$ws.Range("J72").Value = 3198.75
$ws.Range("K72").Value = 30358.5003
$ws.Range("L72").Value = 28788.75
$ws.Range("M72").Value = -25990.5003
$ws.Range("N72").Value = -37524.75
$ws.Range("H98").Value = 844.1539
$ws.Range("I98").Value = 693.2727
$ws.Range("J98").Value = 1674
$ws.Range("K98").Value = 693.2727
$ws.Range("L98").Value = 1674
$ws.Range("M98").Value = 804.7273
$ws.Range("N98").Value = -4670
$ws.Range("H122").Value = 844.1539
$ws.Range("I122").Value = 693.2727
$ws.Range("J122").Value = 1674
$ws.Range("K122").Value = 2079.8181
$ws.Range("L122").Value = 5022
$ws.Range("M122").Value = 370.1819
$ws.Range("N122").Value = -9922
$ws.Range("H132").Value = 3621
$ws.Range("I132").Value = 3414.5
$ws.Range("K132").Value = 10243.5
$ws.Range("M132").Value = -7713.5
$ws.Range("H137").Value = 6897811
$ws.Range("I137").Value = 972.1177
$ws.Range("J137").Value = 16668333
$ws.Range("K137").Value = 2916.3531
$ws.Range("L137").Value = 50004999
$ws.Range("M137").Value = -366.3531000000003
$ws.Range("N137").Value = -50010099
$ws.Range("H141").Value = 2198.75
$ws.Range("I141").Value = 1798.5714
$ws.Range("K141").Value = 5395.7142
$ws.Range("M141").Value = -215.7142000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 33335836
$ws.Range("I61").Value = 45457050
$ws.Range("K61").Value = 45457050
$ws.Range("M61").Value = -45456838
$ws.Range("H63").Value = 50001420
$ws.Range("I63").Value = 62500900
$ws.Range("J63").Value = 3503
$ws.Range("K63").Value = 62500900
$ws.Range("L63").Value = 3503
$ws.Range("M63").Value = -62500214
$ws.Range("N63").Value = -4875
$ws.Range("H66").Value = 50001420
$ws.Range("I66").Value = 62500900
$ws.Range("J66").Value = 3503
$ws.Range("K66").Value = 312504500
$ws.Range("L66").Value = 17515
$ws.Range("M66").Value = -312501068
$ws.Range("N66").Value = -24379
$ws.Range("H74").Value = 25005868
$ws.Range("I74").Value = 38465890
$ws.Range("J74").Value = 8688.429
$ws.Range("K74").Value = 38465890
$ws.Range("L74").Value = 8688.429
$ws.Range("M74").Value = -38465016
$ws.Range("N74").Value = -10436.429
$ws.Range("H77").Value = 25005868
$ws.Range("I77").Value = 38465890
$ws.Range("J77").Value = 8688.429
$ws.Range("K77").Value = 192329450
$ws.Range("L77").Value = 43442.145
$ws.Range("M77").Value = -192325082
$ws.Range("N77").Value = -52178.145
$ws.Range("H122").Value = 4245.2617
$ws.Range("I122").Value = 4651.853
$ws.Range("J122").Value = 2517.25
$ws.Range("K122").Value = 13955.559
$ws.Range("L122").Value = 7551.75
$ws.Range("M122").Value = -11505.559
$ws.Range("N122").Value = -12451.75
$ws.Range("H136").Value = 33335836
$ws.Range("I136").Value = 45457050
$ws.Range("K136").Value = 136371150
$ws.Range("M136").Value = -136368600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10421505
$ws.Range("I31").Value = 5152.926
$ws.Range("J31").Value = 66669810
$ws.Range("K31").Value = 5152.926
$ws.Range("L31").Value = 66669810
$ws.Range("M31").Value = -4857.926
$ws.Range("N31").Value = -66670400
$ws.Range("H34").Value = 10421505
$ws.Range("I34").Value = 5152.926
$ws.Range("J34").Value = 66669810
$ws.Range("K34").Value = 5152.926
$ws.Range("L34").Value = 66669810
$ws.Range("M34").Value = -4950.926
$ws.Range("N34").Value = -66670214
$ws.Range("H107").Value = 444.93332
$ws.Range("I107").Value = 444.93332
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 444.93332
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1475.06668
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 2125.3333
$ws.Range("I122").Value = 2352.4
$ws.Range("K122").Value = 7057.200000000001
$ws.Range("M122").Value = -4607.200000000001
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39.153847
$ws.Range("I12").Value = 23.6
$ws.Range("K12").Value = 70.80000000000001
$ws.Range("M12").Value = 102.2
$ws.Range("H113").Value = 755.6667
$ws.Range("I113").Value = 441.73077
$ws.Range("J113").Value = 1383.5385
$ws.Range("K113").Value = 1325.19231
$ws.Range("L113").Value = 4150.6155
$ws.Range("M113").Value = 844.8076900000001
$ws.Range("N113").Value = -8490.6155
$ws.Range("H125").Value = 3387.5
$ws.Range("I125").Value = 1200
$ws.Range("J125").Value = 4116.6665
$ws.Range("K125").Value = 3600
$ws.Range("L125").Value = 12349.9995
$ws.Range("M125").Value = 1320
$ws.Range("N125").Value = -22189.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 14090.909
$ws.Range("J118").Value = 14090.909
$ws.Range("L118").Value = 14090.909
$ws.Range("N118").Value = -17404.909
$ws.Range("H122").Value = 2223802.8
$ws.Range("I122").Value = 3175699
$ws.Range("J122").Value = 2711.7778
$ws.Range("K122").Value = 9527097
$ws.Range("L122").Value = 8135.3334
$ws.Range("M122").Value = -9524647
$ws.Range("N122").Value = -13035.3334
$ws.Range("H138").Value = 55058.285
$ws.Range("J138").Value = 55058.285
$ws.Range("L138").Value = 55058.285
$ws.Range("N138").Value = -65338.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3104.6177
$ws.Range("I40").Value = 3050.28
$ws.Range("J40").Value = 3255.5557
$ws.Range("K40").Value = 3050.28
$ws.Range("L40").Value = 3255.5557
$ws.Range("M40").Value = -2914.28
$ws.Range("N40").Value = -3527.5557
$ws.Range("H82").Value = 1331.1212
$ws.Range("I82").Value = 940.94446
$ws.Range("J82").Value = 1799.3334
$ws.Range("K82").Value = 940.94446
$ws.Range("L82").Value = 1799.3334
$ws.Range("M82").Value = -579.94446
$ws.Range("N82").Value = -2521.3334
$ws.Range("H85").Value = 1331.1212
$ws.Range("I85").Value = 940.94446
$ws.Range("J85").Value = 1799.3334
$ws.Range("K85").Value = 940.94446
$ws.Range("L85").Value = 1799.3334
$ws.Range("M85").Value = 307.05554
$ws.Range("N85").Value = -4295.3334
$ws.Range("H136").Value = 33335572
$ws.Range("I136").Value = 50002812
$ws.Range("J136").Value = 1093
$ws.Range("K136").Value = 150008436
$ws.Range("L136").Value = 3279
$ws.Range("M136").Value = -150005886
$ws.Range("N136").Value = -8379

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1250.2413
$ws.Range("I132").Value = 1097.8776
$ws.Range("J132").Value = 2079.7778
$ws.Range("K132").Value = 3293.6328
$ws.Range("L132").Value = 6239.3334
$ws.Range("M132").Value = -763.6328000000003
$ws.Range("N132").Value = -11299.3334
$ws.Range("H136").Value = 1150.44
$ws.Range("I136").Value = 1266.0588
$ws.Range("J136").Value = 904.75
$ws.Range("K136").Value = 3798.1764
$ws.Range("L136").Value = 2714.25
$ws.Range("M136").Value = -1248.1764
$ws.Range("N136").Value = -7814.25
$ws.Range("H137").Value = 63746.25
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 63746.25
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 63746.25
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -73946.25
